# Weekly update: insert two new "Papa" price records at rows 415-416,
# pushing the existing rows 415-452 down to 417-454.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 415 (each Insert() pushes the row at
# that index, and everything below it, down by one).
$ws.Rows.Item(415).Insert()
$ws.Rows.Item(415).Insert()

# ---- New row 415 ----
$ws.Cells.Item(415, 1).Value = 8
$ws.Cells.Item(415, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(415, 3).Value = "Coquimbo"
$ws.Cells.Item(415, 4).Value = 44769
$ws.Cells.Item(415, 5).Value = 4
$ws.Cells.Item(415, 6).Value = 100114001
$ws.Cells.Item(415, 7).Value = "Papa"
$ws.Cells.Item(415, 8).Value = "Asterix"
$ws.Cells.Item(415, 9).Value = "1a (cosecha)"
$ws.Cells.Item(415, 10).Value = 2000
$ws.Cells.Item(415, 11).Value = 13000
$ws.Cells.Item(415, 12).Value = 14000
$ws.Cells.Item(415, 13).Value = 13500
$ws.Cells.Item(415, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(415, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(415, 16).Value = 540
$ws.Cells.Item(415, 17).Value = 25
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# ---- New row 416 ----
$ws.Cells.Item(416, 1).Value = 8
$ws.Cells.Item(416, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(416, 3).Value = "Coquimbo"
$ws.Cells.Item(416, 4).Value = 44769
$ws.Cells.Item(416, 5).Value = 4
$ws.Cells.Item(416, 6).Value = 100114001
$ws.Cells.Item(416, 7).Value = "Papa"
$ws.Cells.Item(416, 8).Value = "Rodeo"
$ws.Cells.Item(416, 9).Value = "1a (guarda)"
$ws.Cells.Item(416, 10).Value = 2000
$ws.Cells.Item(416, 11).Value = 11000
$ws.Cells.Item(416, 12).Value = 12000
$ws.Cells.Item(416, 13).Value = 11500
$ws.Cells.Item(416, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(416, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(416, 16).Value = 460
$ws.Cells.Item(416, 17).Value = 25
$ws.Cells.Item(416, 18).Value = "Hortaliza"

# Keep the date cells formatted like the rest of column D.
$ws.Range("D415:D416").NumberFormat = $ws.Range("D417").NumberFormat
